# Updated Code for Sign In, Create Account and Left Hand Panel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 = CreateAccount: Runmode flips from Y to N
$ws.Range("C2").Value = "N"

# Row 8 = LeftHandPanel: Runmode flips from N to Y
$ws.Range("C8").Value = "Y"

# Update the active cell selection to I10
$ws.Range("I10").Select()
